# Refactor code structure for improved readability and maintainability
#
# The "Bilde" (picture path) column previously referenced per-artist image
# files (/album_covers/artist_N.png). These are renamed to per-group image
# files (/album_covers/gruppe_N.png) where N matches the group's row number
# (row 2 = Gruppe 1 = gruppe_1.png, ... row 21 = Gruppe 20 = gruppe_20.png).
# The extra 21st group row is cleared out, leaving 20 groups total.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update the B column (image path) for rows 2..21 to the new gruppe_N.png
# naming scheme, where N = row number - 1 (matches the group number in
# column A, which is left untouched).
for ($row = 2; $row -le 21; $row++) {
    $n = $row - 1
    $ws.Cells.Item($row, 2).Value = "/album_covers/gruppe_$n.png"
}

# The former 21st group (row 22) is removed: clear its Name/Picture/Score
# cell contents while leaving the row's formatting (styles) intact.
$ws.Cells.Item(22, 1).Value = $null
$ws.Cells.Item(22, 2).Value = $null
$ws.Cells.Item(22, 3).Value = $null

# Adjust column widths: Name column widened, Picture column widened
# slightly, and no longer using "best fit" auto-sizing for column A.
$ws.Columns.Item(1).ColumnWidth = 37.41796875
$ws.Columns.Item(2).ColumnWidth = 25.91796875

# Update the active selection to match the saved state (selection on B23
# instead of E22).
$ws.Range("B23").Select()
